# Auto-generated edit script: restores market-price-derived columns (H-N)
# for the Leve Profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the refreshed market data snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 2711.0833
$ws.Cells.Item(15, 9).Value = 2711.0833
$ws.Cells.Item(15, 11).Value = 8133.249899999999
$ws.Cells.Item(15, 13).Value = -7964.249899999999
$ws.Cells.Item(32, 8).Value = 11365856
$ws.Cells.Item(32, 10).Value = 2793.5
$ws.Cells.Item(32, 12).Value = 2793.5
$ws.Cells.Item(32, 14).Value = -3445.5
$ws.Cells.Item(74, 8).Value = 15248
$ws.Cells.Item(74, 9).Value = 15531.2
$ws.Cells.Item(74, 10).Value = 11000
$ws.Cells.Item(74, 11).Value = 15531.2
$ws.Cells.Item(74, 12).Value = 11000
$ws.Cells.Item(74, 13).Value = -14595.2
$ws.Cells.Item(74, 14).Value = -12872
$ws.Cells.Item(76, 8).Value = 3499
$ws.Cells.Item(76, 9).Value = 3499
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 3499
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -3184
$ws.Cells.Item(76, 14).ClearContents()
$ws.Cells.Item(77, 8).Value = 15248
$ws.Cells.Item(77, 9).Value = 15531.2
$ws.Cells.Item(77, 10).Value = 11000
$ws.Cells.Item(77, 11).Value = 77656
$ws.Cells.Item(77, 12).Value = 55000
$ws.Cells.Item(77, 13).Value = -72976
$ws.Cells.Item(77, 14).Value = -64360
$ws.Cells.Item(79, 8).Value = 3499
$ws.Cells.Item(79, 9).Value = 3499
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 3499
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -2407
$ws.Cells.Item(79, 14).ClearContents()
$ws.Cells.Item(123, 8).Value = 110959
$ws.Cells.Item(123, 10).Value = 110959
$ws.Cells.Item(123, 12).Value = 110959
$ws.Cells.Item(123, 14).Value = -120759
$ws.Cells.Item(137, 8).Value = 2216.853
$ws.Cells.Item(137, 9).Value = 2178.2856
$ws.Cells.Item(137, 10).Value = 2396.8333
$ws.Cells.Item(137, 11).Value = 6534.8568
$ws.Cells.Item(137, 12).Value = 7190.499899999999
$ws.Cells.Item(137, 13).Value = -3984.8568
$ws.Cells.Item(137, 14).Value = -12290.4999
$ws.Cells.Item(138, 8).Value = 2661.7
$ws.Cells.Item(138, 9).Value = 910.3077
$ws.Cells.Item(138, 10).Value = 5914.2856
$ws.Cells.Item(138, 11).Value = 2730.9231
$ws.Cells.Item(138, 12).Value = 17742.8568
$ws.Cells.Item(138, 13).Value = 2409.0769
$ws.Cells.Item(138, 14).Value = -28022.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1842760
$ws.Cells.Item(32, 9).Value = 844248.2
$ws.Cells.Item(32, 10).Value = 23810022
$ws.Cells.Item(32, 11).Value = 844248.2
$ws.Cells.Item(32, 12).Value = 23810022
$ws.Cells.Item(32, 13).Value = -843961.2
$ws.Cells.Item(32, 14).Value = -23810596
$ws.Cells.Item(132, 8).Value = 4650.5557
$ws.Cells.Item(132, 9).Value = 7618.5
$ws.Cells.Item(132, 11).Value = 22855.5
$ws.Cells.Item(132, 13).Value = -20325.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(74, 8).Value = 45000
$ws.Cells.Item(74, 10).Value = 45000
$ws.Cells.Item(74, 12).Value = 45000
$ws.Cells.Item(74, 14).Value = -46872
$ws.Cells.Item(77, 8).Value = 45000
$ws.Cells.Item(77, 10).Value = 45000
$ws.Cells.Item(77, 12).Value = 135000
$ws.Cells.Item(77, 14).Value = -144360
$ws.Cells.Item(81, 8).Value = 30780.857
$ws.Cells.Item(81, 10).Value = 33292.832
$ws.Cells.Item(81, 12).Value = 33292.832
$ws.Cells.Item(81, 14).Value = -35414.832
$ws.Cells.Item(84, 8).Value = 30780.857
$ws.Cells.Item(84, 10).Value = 33292.832
$ws.Cells.Item(84, 12).Value = 99878.49600000001
$ws.Cells.Item(84, 14).Value = -110486.496
$ws.Cells.Item(94, 8).Value = 74074540
$ws.Cells.Item(94, 9).Value = 83333730
$ws.Cells.Item(94, 11).Value = 83333730
$ws.Cells.Item(94, 13).Value = -83333279
$ws.Cells.Item(105, 8).Value = 18573152
$ws.Cells.Item(105, 9).Value = 2001425.8
$ws.Cells.Item(105, 10).Value = 27779666
$ws.Cells.Item(105, 11).Value = 2001425.8
$ws.Cells.Item(105, 12).Value = 27779666
$ws.Cells.Item(105, 13).Value = -1999678.8
$ws.Cells.Item(105, 14).Value = -27783160
$ws.Cells.Item(134, 8).Value = 2145.16
$ws.Cells.Item(134, 9).Value = 1314.2307
$ws.Cells.Item(134, 10).Value = 3045.3333
$ws.Cells.Item(134, 11).Value = 3942.6921
$ws.Cells.Item(134, 12).Value = 9135.999899999999
$ws.Cells.Item(134, 13).Value = -1407.6921
$ws.Cells.Item(134, 14).Value = -14205.9999
$ws.Cells.Item(139, 8).Value = 105000
$ws.Cells.Item(139, 10).Value = 105000
$ws.Cells.Item(139, 12).Value = 105000
$ws.Cells.Item(139, 14).Value = -115280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2780640.2
$ws.Cells.Item(31, 9).Value = 2175.6333
$ws.Cells.Item(31, 11).Value = 2175.6333
$ws.Cells.Item(31, 13).Value = -1880.6333
$ws.Cells.Item(34, 8).Value = 2780640.2
$ws.Cells.Item(34, 9).Value = 2175.6333
$ws.Cells.Item(34, 11).Value = 2175.6333
$ws.Cells.Item(34, 13).Value = -1973.6333
$ws.Cells.Item(52, 8).Value = 29960
$ws.Cells.Item(52, 10).Value = 29960
$ws.Cells.Item(52, 12).Value = 29960
$ws.Cells.Item(52, 14).Value = -30548
$ws.Cells.Item(58, 8).Value = 2175.5881
$ws.Cells.Item(58, 10).Value = 2599.4
$ws.Cells.Item(58, 12).Value = 2599.4
$ws.Cells.Item(58, 14).Value = -3005.4
$ws.Cells.Item(86, 8).Value = 6856.1113
$ws.Cells.Item(86, 10).Value = 5953.5
$ws.Cells.Item(86, 12).Value = 5953.5
$ws.Cells.Item(86, 14).Value = -8199.5
$ws.Cells.Item(89, 8).Value = 6856.1113
$ws.Cells.Item(89, 10).Value = 5953.5
$ws.Cells.Item(89, 12).Value = 29767.5
$ws.Cells.Item(89, 14).Value = -40999.5
$ws.Cells.Item(92, 8).Value = 32250
$ws.Cells.Item(92, 10).Value = 32250
$ws.Cells.Item(92, 12).Value = 32250
$ws.Cells.Item(92, 14).Value = -37242
$ws.Cells.Item(99, 8).Value = 5855.2856
$ws.Cells.Item(99, 10).Value = 6198.2
$ws.Cells.Item(99, 12).Value = 6198.2
$ws.Cells.Item(99, 14).Value = -9194.200000000001
$ws.Cells.Item(126, 8).Value = 5855.2856
$ws.Cells.Item(126, 10).Value = 6198.2
$ws.Cells.Item(126, 12).Value = 18594.6
$ws.Cells.Item(126, 14).Value = -23534.6
$ws.Cells.Item(136, 8).Value = 2175.5881
$ws.Cells.Item(136, 10).Value = 2599.4
$ws.Cells.Item(136, 12).Value = 7798.200000000001
$ws.Cells.Item(136, 14).Value = -12898.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 1544.8
$ws.Cells.Item(137, 9).Value = 1544.8
$ws.Cells.Item(137, 11).Value = 4634.4
$ws.Cells.Item(137, 13).Value = 465.6000000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 1299.75
$ws.Cells.Item(2, 9).Value = 3400
$ws.Cells.Item(2, 10).Value = 39.6
$ws.Cells.Item(2, 11).Value = 3400
$ws.Cells.Item(2, 12).Value = 39.6
$ws.Cells.Item(2, 13).Value = -3287
$ws.Cells.Item(2, 14).Value = -265.6
$ws.Cells.Item(80, 8).Value = 142862130
$ws.Cells.Item(80, 9).Value = 333336670
$ws.Cells.Item(80, 11).Value = 333336670
$ws.Cells.Item(80, 13).Value = -333335672
$ws.Cells.Item(83, 8).Value = 142862130
$ws.Cells.Item(83, 9).Value = 333336670
$ws.Cells.Item(83, 11).Value = 1666683350
$ws.Cells.Item(83, 13).Value = -1666678358
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 8320.362999999999
$ws.Cells.Item(126, 9).Value = 2070.75
$ws.Cells.Item(126, 11).Value = 6212.25
$ws.Cells.Item(126, 13).Value = -3742.25
$ws.Cells.Item(132, 8).Value = 2363.3333
$ws.Cells.Item(132, 9).Value = 1894.6666
$ws.Cells.Item(132, 10).Value = 2480.5
$ws.Cells.Item(132, 11).Value = 5683.9998
$ws.Cells.Item(132, 12).Value = 7441.5
$ws.Cells.Item(132, 13).Value = -3153.9998
$ws.Cells.Item(132, 14).Value = -12501.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2370.25
$ws.Cells.Item(46, 9).Value = 2121.2856
$ws.Cells.Item(46, 11).Value = 2121.2856
$ws.Cells.Item(46, 13).Value = -1933.2856
$ws.Cells.Item(93, 8).Value = 3134.7896
$ws.Cells.Item(93, 9).Value = 2982.077
$ws.Cells.Item(93, 10).Value = 3465.6667
$ws.Cells.Item(93, 11).Value = 2982.077
$ws.Cells.Item(93, 12).Value = 3465.6667
$ws.Cells.Item(93, 13).Value = -1734.077
$ws.Cells.Item(93, 14).Value = -5961.6667
$ws.Cells.Item(122, 8).Value = 4932.68
$ws.Cells.Item(122, 9).Value = 3474.1333
$ws.Cells.Item(122, 10).Value = 7120.5
$ws.Cells.Item(122, 11).Value = 10422.3999
$ws.Cells.Item(122, 12).Value = 21361.5
$ws.Cells.Item(122, 13).Value = -7972.3999
$ws.Cells.Item(122, 14).Value = -26261.5
$ws.Cells.Item(127, 8).Value = 66440.5
$ws.Cells.Item(127, 10).Value = 66440.5
$ws.Cells.Item(127, 12).Value = 66440.5
$ws.Cells.Item(127, 14).Value = -76360.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 41124.5
$ws.Cells.Item(86, 10).Value = 41124.5
$ws.Cells.Item(86, 12).Value = 41124.5
$ws.Cells.Item(86, 14).Value = -43370.5
$ws.Cells.Item(89, 8).Value = 41124.5
$ws.Cells.Item(89, 10).Value = 41124.5
$ws.Cells.Item(89, 12).Value = 205622.5
$ws.Cells.Item(89, 14).Value = -216854.5
$ws.Cells.Item(122, 8).Value = 22729772
$ws.Cells.Item(122, 9).Value = 2854.7144
$ws.Cells.Item(122, 10).Value = 62501876
$ws.Cells.Item(122, 11).Value = 8564.143199999999
$ws.Cells.Item(122, 12).Value = 187505628
$ws.Cells.Item(122, 13).Value = -6114.143199999999
$ws.Cells.Item(122, 14).Value = -187510528
$ws.Cells.Item(132, 8).Value = 2939.5806
$ws.Cells.Item(132, 9).Value = 2872.6538
$ws.Cells.Item(132, 10).Value = 3287.6
$ws.Cells.Item(132, 11).Value = 8617.9614
$ws.Cells.Item(132, 12).Value = 9862.799999999999
$ws.Cells.Item(132, 13).Value = -6087.9614
$ws.Cells.Item(132, 14).Value = -14922.8
